# Updates the cryptos list: refreshed Price (D) / Volume(1h) (E) figures for most
# rows, plus a data refresh for rows 50-51 where the coin ordering changed
# (Quant and PaxDollar swapped rank positions), per the upstream GitHub Action run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number -> ordered list of (column, new value) pairs to write.
# D-column (Price) values are prefixed with a literal leading apostrophe so Excel
# keeps them as text exactly as scraped (e.g. "1.00", "170.20") instead of silently
# re-parsing them as numbers and dropping significant trailing/decimal digits.
$rowUpdates = @(
    @{Row=2; Cells=@(@{Col='D';Val='35.630.65'}, @{Col='E';Val='  +3.37%  '})}
    @{Row=3; Cells=@(@{Col='D';Val='1.859.95'}, @{Col='E';Val='  +2.83%  '})}
    @{Row=4; Cells=@(@{Col='E';Val='  +0.34%  '})}
    @{Row=5; Cells=@(@{Col='D';Val='230.68'}, @{Col='E';Val='  +2.35%  '})}
    @{Row=6; Cells=@(@{Col='E';Val='  +3.71%  '})}
    @{Row=7; Cells=@(@{Col='E';Val='  +0.33%  '})}
    @{Row=8; Cells=@(@{Col='D';Val='42.60'}, @{Col='E';Val='  +10.73%  '})}
    @{Row=9; Cells=@(@{Col='D';Val='0.309'}, @{Col='E';Val='  +7.40%  '})}
    @{Row=10; Cells=@(@{Col='D';Val='0.0694'}, @{Col='E';Val='  +3.15%  '})}
    @{Row=11; Cells=@(@{Col='E';Val='  +4.13%  '})}
    @{Row=12; Cells=@(@{Col='D';Val='2.130.88'}, @{Col='E';Val='  +2.84%  '})}
    @{Row=13; Cells=@(@{Col='E';Val='  +4.39%  '})}
    @{Row=14; Cells=@(@{Col='D';Val='1.855.68'}, @{Col='E';Val='  +2.64%  '})}
    @{Row=15; Cells=@(@{Col='D';Val='0.678'}, @{Col='E';Val='  +7.53%  '})}
    @{Row=16; Cells=@(@{Col='D';Val='4.72'}, @{Col='E';Val='  +6.84%  '})}
    @{Row=17; Cells=@(@{Col='D';Val='35.645.99'}, @{Col='E';Val='  +3.43%  '})}
    @{Row=18; Cells=@(@{Col='D';Val='70.34'}, @{Col='E';Val='  +3.20%  '})}
    @{Row=19; Cells=@(@{Col='D';Val='248.74'}, @{Col='E';Val='  +2.36%  '})}
    @{Row=20; Cells=@(@{Col='D';Val='0.0₃0803'}, @{Col='E';Val='  +3.99%  '})}
    @{Row=21; Cells=@(@{Col='D';Val='12.23'}, @{Col='E';Val='  +9.36%  '})}
    @{Row=22; Cells=@(@{Col='D';Val='4.71'}, @{Col='E';Val='  +14.64%  '})}
    @{Row=24; Cells=@(@{Col='D';Val='2.19'}, @{Col='E';Val='  -0.51%  '})}
    @{Row=25; Cells=@(@{Col='D';Val='170.20'}, @{Col='E';Val='  -0.04%  '})}
    @{Row=26; Cells=@(@{Col='D';Val='7.98'}, @{Col='E';Val='  +3.20%  '})}
    @{Row=27; Cells=@(@{Col='D';Val='17.89'}, @{Col='E';Val='  +1.28%  '})}
    @{Row=28; Cells=@(@{Col='E';Val='  +2.19%  '})}
    @{Row=29; Cells=@(@{Col='E';Val='  +16.08%  '})}
    @{Row=30; Cells=@(@{Col='D';Val='1.00'}, @{Col='E';Val='  +0.29%  '})}
    @{Row=31; Cells=@(@{Col='D';Val='3.321.18'}, @{Col='E';Val='  +36.69%  '})}
    @{Row=32; Cells=@(@{Col='D';Val='0.0545'}, @{Col='E';Val='  +5.70%  '})}
    @{Row=33; Cells=@(@{Col='D';Val='4.08'}, @{Col='E';Val='  +5.94%  '})}
    @{Row=34; Cells=@(@{Col='D';Val='3.94'}, @{Col='E';Val='  +4.38%  '})}
    @{Row=35; Cells=@(@{Col='D';Val='1.89'}, @{Col='E';Val='  +3.83%  '})}
    @{Row=36; Cells=@(@{Col='D';Val='100.28'}, @{Col='E';Val='  +22.66%  '})}
    @{Row=37; Cells=@(@{Col='D';Val='0.697'}, @{Col='E';Val='  +8.83%  '})}
    @{Row=38; Cells=@(@{Col='D';Val='1.370.05'}, @{Col='E';Val='  +1.16%  '})}
    @{Row=39; Cells=@(@{Col='D';Val='2.49'}, @{Col='E';Val='  +6.84%  '})}
    @{Row=40; Cells=@(@{Col='E';Val='  +3.01%  '})}
    @{Row=41; Cells=@(@{Col='E';Val='  +4.38%  '})}
    @{Row=42; Cells=@(@{Col='E';Val='  +6.25%  '})}
    @{Row=43; Cells=@(@{Col='D';Val='1.27'}, @{Col='E';Val='  +4.38%  '})}
    @{Row=44; Cells=@(@{Col='D';Val='14.74'}, @{Col='E';Val='  +7.33%  '})}
    @{Row=45; Cells=@(@{Col='E';Val='  +1.17%  '})}
    @{Row=46; Cells=@(@{Col='E';Val='  +1.25%  '})}
    @{Row=47; Cells=@(@{Col='D';Val='6.27'}, @{Col='E';Val='  +8.64%  '})}
    @{Row=48; Cells=@(@{Col='D';Val='0.0522'}, @{Col='E';Val='  +2.25%  '})}
    @{Row=49; Cells=@(@{Col='D';Val='2.028.93'}, @{Col='E';Val='  +2.84%  '})}
    @{Row=50; Cells=@(@{Col='B';Val='PaxDollar'}, @{Col='C';Val='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'}, @{Col='D';Val='1.00'}, @{Col='E';Val='  +0.38%  '})}
    @{Row=51; Cells=@(@{Col='B';Val='Quant'}, @{Col='C';Val='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'}, @{Col='D';Val='104.04'}, @{Col='E';Val='  +1.69%  '})}
)

foreach ($rowUpdate in $rowUpdates) {
    foreach ($cellUpdate in $rowUpdate.Cells) {
        $addr = "$($cellUpdate.Col)$($rowUpdate.Row)"
        if ($cellUpdate.Col -eq "D") {
            # Force text so values like '1.00' or '170.20' are not coerced to numbers.
            $ws.Range($addr).Value = "`'" + $cellUpdate.Val
        } else {
            $ws.Range($addr).Value = $cellUpdate.Val
        }
    }
}
